$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
#    We clone the "2021-Q4" sheet (same column layout / styles) so the
#    new sheet inherits identical sheetPr / styles / pageMargins, then
#    overwrite its contents with the 2022-Q1 fund-holding data.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($totalSheet)

# Re-fetch "总计" - its Index shifted by one after the Copy() call above.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# The template ("2021-Q4") has 12 data rows (rows 2-13); the new sheet
# only needs 7 data rows (rows 2-8), so remove the extra rows.
$newSheet.Range("A9:A13").EntireRow.Delete()

# Fill in the new data. Columns B-G hold text values (even the
# numeric-looking ones), column A and H hold real numbers.
$fundRows = @(
    @("159865", "国泰中证畜牧养殖ETF", "24.32", "99.29", "3.95", "0.9606", 7),
    @("001556", "天弘中证500指数增强A", "41.41", "94.29", "1.69", "0.6998", 10),
    @("001557", "天弘中证500指数增强C", "13.97", "94.29", "1.69", "0.2361", 10),
    @("159867", "鹏华中证畜牧养殖ETF", "5.61", "97.87", "3.89", "0.2182", 7),
    @("516760", "平安中证畜牧养殖ETF", "1.45", "97.82", "3.88", "0.0563", 7),
    @("516670", "招商中证畜牧养殖ETF", "1.06", "98.73", "3.92", "0.0416", 7),
    @("005443", "国金量化多策略灵活配置混合", "0.51", "64.10", "0.73", "0.0037", 7)
)

# Force columns B:G to be stored as text (so numeric-looking strings
# like "24.32" don't get silently converted to numbers), then clear
# the direct formatting again afterwards so no stray style index is
# left behind on the cells.
$textRange = $newSheet.Range("B2:G8")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

$textRange.ClearFormats()

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q1
#    above the existing 2021-Q4 row, and renumber the index column.
# ------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Copy formatting from the row below (the old first data row) so the
# new row matches existing look (style on column A, no style on B:D).
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 7
$totalSheet.Cells.Item(2, 4).Value = 2.22

# Renumber the index column (A) for the rows that shifted down.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
